$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.605.20"
$ws.Range("D3").Value = "3.628.45"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.76"
$ws.Range("E6").Value = "  -4.16%  "
$ws.Range("D7").Value = "3.621.20"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -5.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.79"
$ws.Range("E11").Value = "  +15.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.606"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.43"
$ws.Range("E13").Value = "  -4.49%  "
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").Value = "4.211.11"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "674.47"
$ws.Range("E16").Value = "  -4.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.94"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "3.623.75"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "70.590.03"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.77"
$ws.Range("E21").Value = "  -4.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.47"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.937"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.12"
$ws.Range("E24").Value = "  -4.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.77"
$ws.Range("E25").Value = "  -5.37%  "
$ws.Range("E26").Value = "  -3.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.87"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.64"
$ws.Range("E30").Value = "  -3.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.11"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.62"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.38"
$ws.Range("E34").Value = "  -7.01%  "
$ws.Range("E35").Value = "  -4.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "580.87"
$ws.Range("E36").Value = "  -2.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.08"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.43"
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0454"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").Value = "3.565.87"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("E43").Value = "  -3.27%  "
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.45"
$ws.Range("E45").Value = "  -4.47%  "
$ws.Range("D46").Value = "0.0₃0731"
$ws.Range("E46").Value = "  -6.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  -4.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.76"
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("E51").Value = "  -2.72%  "
